# Update Name of Algo
# Applies updated numeric values (as produced by a re-run of the RandomForest
# imputation algorithm) to the existing result cells on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value  = -11.29639999999999
$ws.Range("B3").Value  = 5.992499999999989
$ws.Range("C5").Value  = -13.831
$ws.Range("E5").Value  = 12.07529999999999
$ws.Range("E9").Value  = 14.26450000000001
$ws.Range("E11").Value = 14.00819999999999
$ws.Range("B14").Value = 9.108900000000004
$ws.Range("B16").Value = 9.922100000000007
$ws.Range("C16").Value = -11.76570000000001
$ws.Range("E17").Value = 13.27140000000001
$ws.Range("B21").Value = 5.471799999999996
$ws.Range("E21").Value = 13.35879999999999
$ws.Range("B23").Value = 5.407500000000002
$ws.Range("B25").Value = 6.031499999999994
